# Update countries & provincias Spain
# Refreshes the COVID country table: updates several countries' case counts,
# re-sorts two countries into their correct alphabetical/ranking position
# (Benin moves up next to Libia/Tunez; Curazao moves up next to Timor
# Oriental), and bumps the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 04:32"

# --- Bolivia (row 43): refreshed case counts ---
$ws.Range("B43").Value = 44113
$ws.Range("C43").Value = 1129
$ws.Range("D43").Value = 13354
$ws.Range("E43").Value = 29121
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 61
$ws.Range("H43").Value = 1638

# --- Australia (row 73): refreshed case counts ---
$ws.Range("B73").Value = 9074
$ws.Range("C73").Value = 15
$ws.Range("D73").Value = 7576
$ws.Range("E73").Value = 1392
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 106

# --- Benin reorder: Benin moves to right after Libia (row 129), pushing
#     Tunez / Suazilandia / Ruanda down one row each (rows 130-133) ---
$ws.Range("A130").Value = "Benin"
$ws.Range("B130").Value = 1285
$ws.Range("C130").Value = 86
$ws.Range("D130").Value = 333
$ws.Range("E130").Value = 929
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 2
$ws.Range("H130").Value = 23

$ws.Range("A131").Value = "Tunez"
$ws.Range("B131").Value = 1231
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 1055
$ws.Range("E131").Value = 126
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 50

$ws.Range("A132").Value = "Suazilandia"
$ws.Range("B132").Value = 1213
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 609
$ws.Range("E132").Value = 587
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 17

$ws.Range("A133").Value = "Ruanda"
$ws.Range("B133").Value = 1210
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 623
$ws.Range("E133").Value = 584
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 3

# Row 134 (Jordania) is unaffected by the Benin reorder - left untouched.

# --- Lesoto (row 179): refreshed case counts ---
$ws.Range("D179").Value = 20
$ws.Range("E179").Value = 113

# --- Curazao reorder: Curazao moves to right after San Vicente y las
#     Granadinas (row 198), swapping places with Timor Oriental
#     (rows 199-200) ---
$ws.Range("A199").Value = "Curazao"
$ws.Range("B199").Value = 25
$ws.Range("C199").Value = 2
$ws.Range("D199").Value = 24
$ws.Range("E199").Value = 0
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1

$ws.Range("A200").Value = "Timor Oriental"
$ws.Range("B200").Value = 24
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 24
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0
